$d = $word.ActiveDocument

# --- 1. Update the report date in the Heading1 paragraph ---------------
[void]$d.Content.Find.Execute("December 04, 2024", $true, $false, $false, $false, $false,
                               $true, 1, $false, "December 11, 2024", 2)

# --- 2. Locate the paragraphs/tables that describe the per-document ----
#        results before we start structurally mutating the body (table
#        deletes leave the Paragraphs collection in a stale state, so we
#        grab every reference we need up front).
$ethPara = $null
$irlPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "ETH-008-lt-leds-2023.pdf*") {
        $ethPara = $p
    } elseif ($t -like "IRL-003-lt-leds-2023.pdf*") {
        $irlPara = $p
    }
}

# --- 3. Turn the "ETH-008-lt-leds-2023.pdf" heading into the new -------
#        single summary line (style Heading2 -> Heading4, new text).
$ethPara.Style = "Heading 4"
$ethPara.Range.Text = "1 documents (0 total pages) processed in 3.96 seconds"

# --- 4. Drop the "IRL-003-lt-leds-2023.pdf" heading paragraph entirely -
$irlPara.Range.Delete()

# --- 5. Drop the two per-document "Variable / Relevant Quotes" tables --
#        (index 3 = the IRL table, index 2 = the ETH table; deleting the
#        higher index first keeps the lower index valid).
$d.Tables.Item(3).Delete()
$d.Tables.Item(2).Delete()

# --- 6. Update the closing summary paragraph's text ---------------------
#        (re-located via Find since object refs above are now stale).
#        Find.Execute collapses the Range it was called on down to the
#        match, so re-use that same Range object to replace the text.
$tail = $d.Content
$ok = $tail.Find.Execute("2 documents (198 total pages) processed in 9.93 seconds",
                          $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok) {
    $tail.Text = "Unable to process the following PDFs: ['C:\Users\WILLIA~1\AppData\Local\Temp\tmpysbn1gdd\doc\8. Twelfth Plan Document.pdf']"
}
